$d = $word.ActiveDocument

function Replace-CellText($table, $row, $col, $old, $new) {
  $cellRng = $table.Cell($row, $col).Range
  $r = $d.Range($cellRng.Start, $cellRng.End)
  $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1) | Out-Null
}

# --- Table of hours worked per person/deliverable (3rd table in the document) ---
$t = $d.Tables.Item(3)

# Luca row: D5 "1:40 h" -> "1:50 h"
Replace-CellText $t 3 6 "1:40 h" "1:50 h"

# Luca row: TOT "68:15 h" -> "68:25 h"
Replace-CellText $t 3 7 "68:15 h" "68:25 h"

# Zakaria row: D4 "29:30 h" -> "34:30 h"
Replace-CellText $t 4 5 "29:30 h" "34:30 h"

# Zakaria row: TOT "59:45 h" -> "64:45 h"
Replace-CellText $t 4 7 "59:45 h" "64:45 h"

# TOT row: D4 "93:10 h" -> "98:10 h"
Replace-CellText $t 5 5 "93:10 h" "98:10 h"

# TOT row: D5 "3:40 h" -> "3:50 h"
Replace-CellText $t 5 6 "3:40 h" "3:50 h"

# TOT row: TOT "186:50 h" -> "192:00 h"
Replace-CellText $t 5 7 "186:50 h" "192:00 h"

# --- Narrative paragraph about hour distribution (D4 section) ---
# Done in two pieces (split at the same point as the existing "_GoBack" bookmark)
# so the bookmark that sits between "gr" and "andi" is left untouched.

# Part before the bookmark: "...non ci sono state gr" -> "...non ci sono state grandi dispar"
$d.Content.Find.Execute("non ci sono state gr", $true, $false, $false, $false, $false, $true, 0, $false, "non ci sono state grandi dispar", 1) | Out-Null

# Part after the bookmark: "andi disparità ... Come già detto precedentemente" -> "ità a livello di ore e nel D4 ... Come già detto precedentemente"
$oldAfter = "andi disparità a livello di ore e nel D4 Luca si ritrova ad avere un numero maggiore di ore poiché ha realizzato la parte di back-end insieme a Zakaria, oltre a creare la parte degli User Flow, il quale all" + [char]8217 + "inizio era il suo compito. Come già detto precedentemente"
$newAfter = "ità a livello di ore e nel D4 sia Luca che Zakaria si ritrovano ad avere un numero maggiore di ore poiché impegnati sia nella realizzazione di diagrammi (User flow per Luca, i vari diagrammi delle API per Zakaria) che nella realizzazione del back-end. Come già detto precedentemente"
$d.Content.Find.Execute($oldAfter, $true, $false, $false, $false, $false, $true, 0, $false, $newAfter, 1) | Out-Null
